# Applies the "updates to water balance and glm" change to the CHART sheet:
#  - adds a new row 38 (S38:BA38 = 1..35, a helper index row used by the
#    downstream bathymetry/GLM chart calculations)
#  - updates the sheet's view state (active selection / scroll position)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHART")
$ws.Activate()

# New row 38: sequential integers 1..35 across columns S:BA
$row38 = New-Object 'object[,]' 1,35
for ($i = 0; $i -lt 35; $i++) {
    $row38[0, $i] = $i + 1
}
$ws.Range("S38:BA38").Value = $row38

# Update the view: scroll the window so column AC is the left-most visible
# column, then select AY43 (matches the saved sheetView/selection state).
$excel.ActiveWindow.ScrollColumn = 29
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AY43").Select()
